$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; existing rows 14-38 shift down to 15-39.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value = 6
$ws.Cells.Item(14, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(14, 3).Value = "Metropolitana"
$ws.Cells.Item(14, 4).Value = 44791
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = 100112035
$ws.Cells.Item(14, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 16000
$ws.Cells.Item(14, 12).Value = 18000
$ws.Cells.Item(14, 13).Value = 17133
$ws.Cells.Item(14, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(14, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(14, 16).Value = 1142
$ws.Cells.Item(14, 17).Value = 15
$ws.Cells.Item(14, 18).Value = "Hortaliza"
